# Append a new log row (row 31) to the Nalco run-log sheet, mirroring the
# formatting of the existing data rows (e.g. row 30) and filling in the
# new run's details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 30
$newRow = 31

# Clone the formatting (style) of the previous row onto the new row so the
# appended row matches the look of the rest of the log (center/center
# alignment), without creating a brand-new, unused style entry.
$ws.Range("A$($lastRow):H$($lastRow)").Copy()
$ws.Range("A$($newRow):H$($newRow)").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A$newRow").Value = "2025-08-19 06:49:33 UTC"
$ws.Range("B$newRow").Value = "2025-08-19 12:19:33 IST"
$ws.Range("C$newRow").Value = "SKIPPED"
$ws.Range("D$newRow").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E$newRow").Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Range("G$newRow").Value = 0

Write-Host "Appended row $newRow to $($ws.Name)"
